# Apply the EngineBugAndRequest.xlsx update:
# - appends 6 new ISSUE rows (rows 5-10) with their STATUS ("DONE") for row 9
# - moves the active-cell selection to A7 (matches the author's final cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Separate Weapon Animation from Hero Animation"
$ws.Range("A6").Value = "Discrete animation for heroes"
$ws.Range("A7").Value = "Quick Animator portraits are still needed, but should be their own output so that each import doesn" + [char]0x2019 + "t require re-sizing"
$ws.Range("A8").Value = "Decide if spells should be split out in the animations"
$ws.Range("A9").Value = "Level up/level down/reset sprite items or menu debug menu. "
$ws.Range("D9").Value = "DONE"
$ws.Range("A10").Value = "Select battle text file in development mode"

$ws.Range("A7").Select() | Out-Null
